$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A87").Value = "GRT-USD"
